$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of row 40 down into row 41 (same pattern Excel uses
# when a user duplicates the last data row to add a new entry), then
# overwrite the text for the new "DIFF_SEGNALAZIONE_OVER" variable.
$ws.Range("A40:F40").Copy() | Out-Null
$ws.Range("A41:F41").PasteSpecial() | Out-Null

$ws.Range("A41").Value = "CREATE/MODIFY"
$ws.Range("B41").Value = "DIFF_SEGNALAZIONE_OVER"
$ws.Range("C41").Value = "DIFF_SEGNALAZIONE_OVER"
$ws.Range("D41").Value = ""
$ws.Range("E41").Value = "CUSTOMER"
$ws.Range("F41").Value = "DIFF_SEGNALAZIONE_OVER"

$ws.Range("F41").Select() | Out-Null
